# ItemList.xlsx update:
#  - add a new "Key" item category (header in Z1) with its first entry
#    "Iron Key" (id 7001) in Z3/AA3
#  - add a new "Fireball" spell entry (id 4004) to the existing Spell
#    column (Q/R), right below the other spells
#  - leave the final selection on the newly added spell cell (R6), as
#    was left selected by the author after editing

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Key" category header
$ws.Range("Z1").Value = "Key"

# First item of the new "Key" category
$ws.Range("Z3").Value = "Iron Key"
$ws.Range("AA3").Value = 7001

# New spell added under the existing "Spell" category
$ws.Range("Q6").Value = "Fireball"
$ws.Range("R6").Value = 4004

# Match the final selection left in the workbook
$ws.Range("R6").Select()
